$d = $word.ActiveDocument

$ok = $d.Content.Find.Execute("Contoso CipherGuard Sentinel X7 est un produit de sécurité avancé et résilient, méticuleusement conçu pour renforcer l’infrastructure du réseau informatique contre un large éventail de menaces et de vulnérabilités.", $true, $false, $false, $false, $false, $true, 1, $false, "Le Contoso CipherGuard Sentinel X7 est un produit de sécurité avancé et résilient, spécialement conçu pour renforcer l’infrastructure réseau informatique face à un éventail de menaces et de vulnérabilités.", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #0: 'Contoso CipherGuard Sentinel X7 est un p'" }
$ok = $d.Content.Find.Execute("Protection du pare-feu :", $true, $false, $false, $false, $false, $true, 1, $false, "Protection par pare-feu :", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #1: 'Protection du pare-feu :'" }
$ok = $d.Content.Find.Execute(" l’utilisation d’un pare-feu d’inspection avec état, Contoso CipherGuard Sentinel X7 utilise des techniques d’inspection approfondies des paquets.", $true, $false, $false, $false, $false, $true, 1, $false, " s’appuyant sur un pare-feu d’inspection avec état, le Contoso CipherGuard Sentinel X7 a recours à des techniques d’inspection approfondies des paquets.", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #2: ' l’utilisation d’un pare-feu d’inspectio'" }
$ok = $d.Content.Find.Execute("alimenté par des algorithmes d’apprentissage automatique, notre IDPS surveille en permanence les modèles et anomalies du trafic réseau.", $true, $false, $false, $false, $false, $true, 1, $false, "utilisant des algorithmes d’apprentissage automatique, notre IDPS surveille en permanence les modèles et anomalies du trafic réseau.", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #3: 'alimenté par des algorithmes d’apprentis'" }
$ok = $d.Content.Find.Execute("Prise en charge du réseau privé virtuel (VPN) : ", $true, $false, $false, $false, $false, $true, 1, $false, "Prise en charge de réseaux privés virtuels (VPN) : ", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #4: 'Prise en charge du réseau privé virtuel '" }
$ok = $d.Content.Find.Execute("Contoso CipherGuard Sentinel X7 prend en charge les protocoles VPN standard tels que IPsec et OpenVPN.", $true, $false, $false, $false, $false, $true, 1, $false, "le Contoso CipherGuard Sentinel X7 prend en charge les protocoles VPN standard tels que IPsec et OpenVPN.", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #5: 'Contoso CipherGuard Sentinel X7 prend en'" }
$ok = $d.Content.Find.Execute(" utilisation d’une approche de défense multicouche, notre module de sécurité de point de terminaison intègre des fonctionnalités antivirus, anti-programme malveillant et de prévention des intrusions basées sur l’hôte.", $true, $false, $false, $false, $false, $true, 1, $false, " se basant sur une approche de défense multicouche, notre module de sécurité des points de terminaison intègre des fonctionnalités d’antivirus, d’anti-programme malveillant et de prévention des intrusions basées sur l’hôte.", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #6: ' utilisation d’une approche de défense m'" }
$ok = $d.Content.Find.Execute("Authentification utilisateur et contrôle", $true, $false, $false, $false, $false, $true, 1, $false, "Authentification utilisateur et contrôle d’accès", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #7: 'Authentification utilisateur et contrôle'" }
$ok = $d.Content.Find.Execute(" d’accès : Contoso CipherGuard Sentinel X7 prend en charge les mécanismes d’authentification multifacteur (MFA), notamment l’authentification biométrique et l’intégration de cartes à puce.", $true, $false, $false, $false, $false, $true, 1, $false, " : le Contoso CipherGuard Sentinel X7 prend en charge les mécanismes d’authentification multifacteur (MFA), notamment l’authentification biométrique et l’intégration de cartes à puce.", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #8: ' d’accès : Contoso CipherGuard Sentinel '" }
$ok = $d.Content.Find.Execute(" Quad-core 2,5 GHz ou supérieur avec prise en charge de l’accélération matérielle", $true, $false, $false, $false, $false, $true, 1, $false, " quad-core 2,5 GHz ou supérieur avec prise en charge de l’accélération matérielle", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #9: ' Quad-core 2,5 GHz ou supérieur avec pri'" }
$ok = $d.Content.Find.Execute(" 16 Go minimum, ECC (Code de correction d’erreur) recommandé", $true, $false, $false, $false, $false, $true, 1, $false, " 16 Go minimum, ECC (code de correction d’erreur) recommandé", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #10: ' 16 Go minimum, ECC (Code de correction '" }
$ok = $d.Content.Find.Execute(" Double Ethernet Ethernet avec prise en charge des trames jumbo", $true, $false, $false, $false, $false, $true, 1, $false, " deux ports Gigabit Ethernet avec prise en charge des trames Jumbo", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #11: ' Double Ethernet Ethernet avec prise en '" }
$ok = $d.Content.Find.Execute(" TCP/IP, UDP, ICMP, prise en charge IPv6", $true, $false, $false, $false, $false, $true, 1, $false, " prise en charge de TCP/IP, UDP, ICMP, IPv6", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #12: ' TCP/IP, UDP, ICMP, prise en charge IPv6'" }
$ok = $d.Content.Find.Execute(" Interopérabilité avec Cisco, Juniper et d’autres fournisseurs de réseaux majeurs", $true, $false, $false, $false, $false, $true, 1, $false, " interopérabilité avec Cisco, Juniper et d’autres fournisseurs de réseaux majeurs", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #13: ' Interopérabilité avec Cisco, Juniper et'" }
$ok = $d.Content.Find.Execute(" effectuez une évaluation complète des vulnérabilités réseau, notamment les tests d’intrusion et l’analyse des risques.", $true, $false, $false, $false, $false, $true, 1, $false, " effectuez une évaluation complète des vulnérabilités réseau, notamment des tests d’intrusion et une analyse des risques.", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #14: ' effectuez une évaluation complète des v'" }
$ok = $d.Content.Find.Execute(" Déployez Contoso CipherGuard Sentinel X7 sur des serveurs dédiés ou des machines virtuelles, ce qui garantit une utilisation optimale du matériel et l’allocation des ressources.", $true, $false, $false, $false, $false, $true, 1, $false, " déployez le Contoso CipherGuard Sentinel X7 sur des machines virtuelles ou des serveurs dédiés, afin de garantir une utilisation optimale du matériel et une allocation efficace des ressources.", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #15: ' Déployez Contoso CipherGuard Sentinel X'" }
$ok = $d.Content.Find.Execute(" Personnalisez les stratégies de sécurité, les contrôles d’accès et les règles de pare-feu en fonction des exigences organisationnelles.", $true, $false, $false, $false, $false, $true, 1, $false, " personnalisez les stratégies de sécurité, les contrôles d’accès et les règles de pare-feu en fonction des exigences organisationnelles.", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #16: ' Personnalisez les stratégies de sécurit'" }
$ok = $d.Content.Find.Execute(" exécutez un plan de test approfondi, y compris des scénarios d’attaque simulé et des tests de charge, pour valider l’efficacité et les performances de la solution.", $true, $false, $false, $false, $false, $true, 1, $false, " exécutez un plan de test approfondi, notamment des scénarios de simulation d’attaque et des tests de charge, pour valider l’efficacité et les performances de la solution.", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #17: ' exécutez un plan de test approfondi, y '" }
$ok = $d.Content.Find.Execute(" fournissez des sessions de formation approfondies au personnel informatique, couvrant les opérations quotidiennes, les procédures de réponse aux incidents et les tâches de maintenance.", $true, $false, $false, $false, $false, $true, 1, $false, " proposez des sessions de formation approfondie au personnel informatique, couvrant les opérations quotidiennes, les procédures de réponse aux incidents et les tâches de maintenance.", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #18: ' fournissez des sessions de formation ap'" }
$ok = $d.Content.Find.Execute(" Contoso garantit des mises à jour continues du produit, en intégrant les dernières améliorations du renseignement sur les menaces et de la sécurité.", $true, $false, $false, $false, $false, $true, 1, $false, " Contoso garantit des mises à jour continues du produit, intégrant les dernières améliorations du renseignement sur les menaces et de la sécurité.", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #19: ' Contoso garantit des mises à jour conti'" }
$ok = $d.Content.Find.Execute(" Contoso fournit une équipe de support technique dédiée 24/7 pour garantir une assistance rapide pour tout problème technique ou demande de renseignements liés à Contoso CipherGuard Sentinel X7.", $true, $false, $false, $false, $false, $true, 1, $false, " Contoso fournit une équipe de support technique dédiée 24/7 pour garantir une assistance rapide en cas de problème technique ou de demande de renseignements liés au Contoso CipherGuard Sentinel X7.", 2)
if (-not $ok) { Write-Host "REPLACE FAILED #20: ' Contoso fournit une équipe de support t'" }

$rng = $d.Content
$found = $rng.Find.Execute("Système de détection et de prévention des intrusions (IDPS) : ")
if ($found) {
    $rng.Collapse(1)
    $rng.InsertBefore("u")
    $rng.Style = "Default Paragraph Font"
    $rng.Font.Name = "Calibri"
    $rng.Font.NameFarEast = "Calibri"
    $rng.Font.NameBi = "Times New Roman"
    $rng.Font.Bold = $false
    $rng.Font.BoldBi = $true
    $rng.Font.Italic = $false
    $rng.Font.ItalicBi = $false
    $rng.Font.AllCaps = $false
    $rng.Font.SmallCaps = $false
    $rng.Font.StrikeThrough = $false
    $rng.Font.DoubleStrikeThrough = $false
    $rng.Font.Outline = $false
    $rng.Font.Shadow = $false
    $rng.Font.Emboss = $false
    $rng.Font.Engrave = $false
    $rng.Font.Hidden = $false
    $rng.Font.Size = 11
    $rng.Font.Underline = 0
    $rng.Font.Position = 0
    $rng.Font.Spacing = 0
    $rng.Font.Kerning = 0
} else {
    Write-Host "IDPS paragraph not found for u-run insertion"
}

Write-Host "Done."
